$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2), pushing all existing
# rows (and the last row, 49 -> 50) down by one.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# New day's date is one day after the (old) most-recent date, which is
# now in row 3.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-01-08"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
